# Insert a new data row at row 122 (pushing existing rows 122:197 down to 123:198),
# then populate the new row with its price-report values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(122).Insert()

$ws.Range("A122").Value = 4
$ws.Range("B122").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C122").Value = "Los Lagos"
$ws.Range("D122").Value = 44488
$ws.Range("E122").Value = 10
$ws.Range("F122").Value = 100114014
$ws.Range("G122").Value = "Betarraga"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 1000
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = 1100
$ws.Range("N122").Value = "$/paquete 5 unidades"
$ws.Range("O122").Value = "Región del Maule"
$ws.Range("P122").Value = 220
$ws.Range("Q122").Value = 5
$ws.Range("R122").Value = "Hortaliza"
